$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 43663
$ws.Range("A7").NumberFormat = "yyyy/m/d"
$ws.Range("A7").HorizontalAlignment = -4108
$ws.Range("A7").VerticalAlignment = -4108

$ws.Range("B7").Value = "上午08:00-11:30"
$ws.Range("C7").Value = "初步完成zutnlp-entity后端"
$ws.Range("D7").Value = "完成"

$ws.Range("B7:D7").HorizontalAlignment = -4108
$ws.Range("B7:D7").VerticalAlignment = -4108

$ws.Range("D7").Select()
